# Daily update at 8 AM UTC
# Adds the next day's row (row 27) to the Wins Over Time sheet and moves the
# "last row" date-only formatting down from row 26 to the newly appended row 27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26 is no longer the last row, so it reverts to the regular datetime
# number format used by the rest of the data rows (style index 2).
$ws.Range("A26").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row 27 with the next day's data.
$ws.Range("A27").Value2 = 45976
$ws.Range("B27").Value2 = 58
$ws.Range("C27").Value2 = 65
$ws.Range("D27").Value2 = 68

# Row 27 is now the last row, so it gets the date-only number format
# (style index 3) that row 26 previously had.
$ws.Range("A27").NumberFormat = "YYYY-MM-DD"
